$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.460.91'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.55%  '

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.694.31'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.58%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  +0.62%  '

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.50%  '

# Row 6 - XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5550'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +8.80%  '

# Row 7 - USDC
$ws.Range("E7").Value = '  +0.61%  '

# Row 8 - Cardano
$ws.Range("E8").Value = '  +1.82%  '

# Row 9 - Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06486'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.45%  '

# Row 10 - Solana
$ws.Range("E10").Value = '  +1.23%  '

# Row 11 - TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07638'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.56%  '

# Row 12 - Polkadot
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.567'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.21%  '

# Row 13 - WrappedEther
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.683.61'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.84%  '

# Row 14 - Polygon
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5832'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.41%  '

# Row 15 - ShibaInu
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008468'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.34%  '

# Row 16 - Litecoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.35'
$ws.Range("D16").Style = "Normal"

# Row 17 - WrappedBTC
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.523.28'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.45%  '

# Row 18 - Uniswap
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.966'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.88%  '

# Row 19 - Dai
$ws.Range("E19").Value = '  +0.64%  '

# Row 20 - Avalanche
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.00'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.86%  '

# Row 21 - BitcoinCash
$ws.Range("E21").Value = '  +0.64%  '

# Row 22 - Chainlink
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.264'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.19%  '

# Row 23 - BinanceUSD
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.011'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.59%  '

# Row 24 - Monero
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '150.59'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.76%  '

# Row 25 - Stellar
$ws.Range("E25").Value = '  +7.80%  '

# Row 26 - Cosmos
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.925'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.02%  '

# Row 27 - EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.78'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.77%  '

# Row 28 - Toncoin
$ws.Range("E28").Value = '  +7.19%  '

# Row 29 - Hedera
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.06348'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.03%  '

# Row 30 - PancakeSwap
$ws.Range("E30").Value = '  +1.29%  '

# Row 31 - Filecoin
$ws.Range("E31").Value = '  +1.22%  '

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.597'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.43%  '

# Row 33 - LidoDAOToken
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.683'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.36%  '

# Row 34 - ARBITRUM
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.046'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.93%  '

# Row 35 - ImmutableX
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6238'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.16%  '

# Row 36 - HuobiToken
$ws.Range("E36").Value = '  +1.48%  '

# Row 37 - MXToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.725'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.55%  '

# Row 38 - FraxShare
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.238'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.41%  '

# Row 39 - Maker->VeChain
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01646'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.36%  '

# Row 40 - VeChain->Maker
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.123.53'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.35%  '

# Row 41 - TrustWalletToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8850'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.81%  '

# Row 43 - Quant
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.79'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.61%  '

# Row 44 - RocketPoolETH
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.843.14'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.51%  '

# Row 45 - BabyDogeCoin
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000110'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.99%  '

# Row 46 - Aave
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '57.60'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.33%  '

# Row 47 - EnergySwap
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.236'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.48%  '

# Row 48 - Frax
$ws.Range("E48").Value = '  +0.04%  '

# Row 49 - Cronos
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05285'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.02%  '

# Row 50 - Mantle
$ws.Range("E50").Value = '  +0.58%  '

# Row 51 - Aptos
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.085'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.58%  '
